
$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# =============================================================================
# 1. Update the "总计" (summary) sheet: insert a new 2022-Q4 row at the top
#    of the data (row 2), shifting the existing 2022-Q3 / 2022-Q2 / 2022-Q1
#    rows down by one.
# =============================================================================

# Capture the existing quarter/count/value rows (old rows 2-4: Q3, Q2, Q1).
$oldBD = $summary.Range("B2:D4").Value2

# Shift them down into rows 3-5.
$summary.Range("B3:D5").Value2 = $oldBD

# The new row 5 needs the same formatting (bold/border/centered) on column A
# as the other rows already have; copy it down from row 4.
$summary.Range("A4").Copy()
$summary.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-number the column A sequence (0,1,2,3) for the now 4 data rows.
$summary.Range("A2").Value2 = 0
$summary.Range("A3").Value2 = 1
$summary.Range("A4").Value2 = 2
$summary.Range("A5").Value2 = 3

# Write the new 2022-Q4 figures into row 2.
$summary.Range("B2").Value2 = "2022-Q4"
$summary.Range("C2").Value2 = 34
$summary.Range("D2").Value2 = 17.74

# =============================================================================
# 2. Insert a new worksheet "2022-Q4" right after "总计" (so it becomes the
#    second sheet, shifting the old 2022-Q3 / 2022-Q2 / 2022-Q1 sheets back).
# =============================================================================
$q4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summary)
$q4.Name = "2022-Q4"

$n = 34
$colA = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33)
$colB = @("519702","001645","014038","020001","013430","009804","010488","011568","160212","481010","011251","014039","160215","011321","011252","160610","005313","005314","013627","013626","160603","011331","002259","001428","009246","010489","011569","166007","001789","011332","001884","011907","487016","016616")
$colC = @("交银趋势优先混合A","国泰大健康股票A","交银启诚混合A","国泰金鹰增长灵活配置混合","交银趋势优先混合C","国泰研究优势混合","鹏华优选成长混合A","鹏华产业升级混合A","国泰估值优势混合（LOF）A","工银中小盘混合","华安聚嘉精选混合A","交银启诚混合C","国泰价值经典灵活配置混合（LOF）","国泰大健康股票C","华安聚嘉精选混合C","鹏华动力增长混合（LOF）","万家中证1000指数增强A","万家中证1000指数增强C","华夏周期驱动混合C","华夏周期驱动混合A","鹏华普天收益混合","鹏华远见成长混合A","鹏华健康环保灵活配置混合","工银灵活配置混合B","大摩ESG量化混合","鹏华优选成长混合C","鹏华产业升级混合C","中欧互通精选混合A","国泰量化收益灵活配置混合A","鹏华远见成长混合C","中欧互通精选混合E","国泰量化收益灵活配置混合C","工银灵活配置混合A","国泰估值优势混合（LOF）C")
$colD = @("99.56","36.92","24.58","20.52","20.19","15.37","29.13","16.84","9.03","15.71","16.85","7.22","5.95","5.56","11.75","12.89","22.07","19.61","4.87","4.12","4.20","1.72","1.93","2.65","2.39","0.77","0.27","0.57","0.47","0.18","0.03","0.01","0.00","0.00")
$colE = @("81.36","94.51","81.04","94.14","81.36","93.94","63.66","68.55","93.64","91.30","84.72","81.04","94.08","94.51","84.72","53.20","94.13","94.13","86.45","86.45","70.03","69.04","76.75","71.03","92.14","63.66","68.55","91.42","85.24","69.04","91.42","85.24","71.03","93.64")
$colF = @("5.13","6.39","5.51","6.02","5.13","5.87","2.42","3.87","5.91","3.33","2.93","5.51","6.07","6.39","2.93","2.19","1.08","1.08","3.52","3.52","2.36","3.04","2.68","1.71","1.12","2.42","3.87","1.65","1.82","3.04","1.65","1.82","1.71","5.91")
$colG = @("5.1074","2.3592","1.3544","1.2353","1.0357","0.9022","0.7049","0.6517","0.5337","0.5231","0.4937","0.3978","0.3612","0.3553","0.3443","0.2823","0.2384","0.2118","0.1714","0.1450","0.0991","0.0523","0.0517","0.0453","0.0268","0.0186","0.0104","0.0094","0.0086","0.0055","0.0005","0.0002","0","0")
$colH = @(1,9,1,9,1,9,10,3,5,5,3,1,9,9,3,9,2,2,6,6,10,4,10,5,10,10,3,8,7,4,8,7,5,5)


# Header row (row 1) of the new sheet.
$headers = New-Object 'object[,]' 1,7
$headers[0,0] = "基金代码"
$headers[0,1] = "基金名称"
$headers[0,2] = "基金规模"
$headers[0,3] = "股票总仓位"
$headers[0,4] = "仓位占比"
$headers[0,5] = "持有市值(亿元)"
$headers[0,6] = "仓位排名"
$q4.Range("B1:H1").Value2 = $headers

$lastRow = $n + 1

# Column A - sequential index, numeric.
$arrA = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) { $arrA[$i,0] = $colA[$i] }
$q4.Range("A2:A$lastRow").Value2 = $arrA

# Columns B..F - always stored as text (fund codes / names / numeric-looking
# strings whose formatting -- e.g. leading/trailing zeros -- must be kept).
$q4.Range("B2:F$lastRow").NumberFormat = "@"
$arrBF = New-Object 'object[,]' $n,5
for ($i = 0; $i -lt $n; $i++) {
    $arrBF[$i,0] = $colB[$i]
    $arrBF[$i,1] = $colC[$i]
    $arrBF[$i,2] = $colD[$i]
    $arrBF[$i,3] = $colE[$i]
    $arrBF[$i,4] = $colF[$i]
}
$q4.Range("B2:F$lastRow").Value2 = $arrBF

# Column G - held market value. Stored as text EXCEPT when the value is
# exactly zero, in which case it is a genuine number 0.
$arrG = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    if ($colG[$i] -eq "0") {
        $arrG[$i,0] = 0
    } else {
        $arrG[$i,0] = $colG[$i]
    }
}
for ($i = 0; $i -lt $n; $i++) {
    $rowNum = $i + 2
    if ($colG[$i] -ne "0") {
        $q4.Range("G$rowNum").NumberFormat = "@"
    }
}
$q4.Range("G2:G$lastRow").Value2 = $arrG

# Column H - position ranking, numeric.
$arrH = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) { $arrH[$i,0] = $colH[$i] }
$q4.Range("H2:H$lastRow").Value2 = $arrH

# Formatting: header row + column A get the bold/bordered/centered style
# already used on the other quarter sheets (copy it from the summary sheet,
# which still carries the original style definitions).
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$summary.Range("A2").Copy()
$q4.Range("A2:A$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q4.Range("A1").Select()
